# Add the beginner's guide (新手教程) part entries and their dialog lines,
# then leave the selection/navigation state the way the author left it.

$wb = $excel.ActiveWorkbook

$wsPart   = $wb.Worksheets.Item("part")
$wsDialog = $wb.Worksheets.Item("dialog")

# --- part sheet: four new tutorial "part" headings (ids 19-22) ------------
# --- dialog sheet: the dialogue lines that belong to those parts ----------
# The two sheets are interleaved here so that the new entries land in
# xl/sharedStrings.xml in the same order as the reference edit.

$wsPart.Range("A20").Value = 19
$wsPart.Range("B20").Value = "农田新手教程"

$wsDialog.Range("A69").Value = 68
$wsDialog.Range("B69").Value = 19
$wsDialog.Range("C69").Value = 7
$wsDialog.Range("D69").Value = "normal"
$wsDialog.Range("E69").Value = "鼠标左键点击右上角背包图标打开物品栏"

$wsDialog.Range("A70").Value = 69
$wsDialog.Range("B70").Value = 19
$wsDialog.Range("C70").Value = 7
$wsDialog.Range("D70").Value = "normal"
$wsDialog.Range("E70").Value = "左键点击物品栏中的种子"

$wsDialog.Range("A71").Value = 70
$wsDialog.Range("B71").Value = 19
$wsDialog.Range("C71").Value = 7
$wsDialog.Range("D71").Value = "normal"
$wsDialog.Range("E71").Value = "右键点击要种下的地块（只有左边第一列可种，其他田地需要花钱开垦）"

$wsDialog.Range("A72").Value = 71
$wsDialog.Range("B72").Value = 19
$wsDialog.Range("C72").Value = 7
$wsDialog.Range("D72").Value = "normal"
$wsDialog.Range("E72").Value = "种植成功！随时间推移，种植的作物将会成熟，收割后可以出售，是非常重要的经济来源哦~"

$wsPart.Range("A21").Value = 20
$wsPart.Range("B21").Value = "去集市小游戏教程"

$wsDialog.Range("A73").Value = 72
$wsDialog.Range("B73").Value = 20
$wsDialog.Range("C73").Value = 7
$wsDialog.Range("D73").Value = "normal"
$wsDialog.Range("E73").Value = "去集市路途遥远，需要玩一个小小的游戏，AD左右移动，W键跳跃，收集金币可转化成你的资产哦！"

$wsPart.Range("A22").Value = 21
$wsPart.Range("B22").Value = "在集市情绪提示"

$wsDialog.Range("A74").Value = 73
$wsDialog.Range("B74").Value = 21
$wsDialog.Range("C74").Value = 7
$wsDialog.Range("D74").Value = "normal"
$wsDialog.Range("E74").Value = "右边是情绪条，你现在的情绪值有点低，可以通过购买喜欢的物品来提高情绪值。情绪值若太低，有一定几率做出不理智事件哦！"

$wsPart.Range("A23").Value = 22
$wsPart.Range("B23").Value = "正式开始的说明"

$wsDialog.Range("A75").Value = 74
$wsDialog.Range("B75").Value = 22
$wsDialog.Range("C75").Value = 7
$wsDialog.Range("D75").Value = "normal"
$wsDialog.Range("E75").Value = "每天至少要上交家里1块钱 剩余的钱存起来当做学费 需要在九月开学之前攒够学费"

# --- leave the cursor/viewport the way the author left them ---------------
$wsPart.Range("A24").Select()
$wsDialog.Range("E78").Select()
